$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to make room for the header
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "Vorname"
$ws.Range("B1").Value = "Nachname"
$ws.Range("C1").Value = "Karte"

# Autofit column B (Nachname) so it best-fits its contents
$ws.Columns.Item(2).AutoFit()

# Set the active selection to A2
$ws.Range("A2").Select()
